$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7697707.5
$ws.Range("I74").Value = 9095291
$ws.Range("K74").Value = 9095291
$ws.Range("M74").Value = -9094355
$ws.Range("H77").Value = 7697707.5
$ws.Range("I77").Value = 9095291
$ws.Range("K77").Value = 45476455
$ws.Range("M77").Value = -45471775
$ws.Range("H112").Value = 1630.337
$ws.Range("J112").Value = 1638.6364
$ws.Range("L112").Value = 4915.9092
$ws.Range("N112").Value = -7131.9092
$ws.Range("H129").Value = 897.0714
$ws.Range("I129").Value = 398.8
$ws.Range("J129").Value = 964.4054
$ws.Range("K129").Value = 1196.4
$ws.Range("L129").Value = 2893.2162
$ws.Range("M129").Value = 3803.6
$ws.Range("N129").Value = -12893.2162
$ws.Range("H132").Value = 55781868
$ws.Range("I132").Value = 66935976
$ws.Range("J132").Value = 11333.333
$ws.Range("K132").Value = 200807928
$ws.Range("L132").Value = 33999.999
$ws.Range("M132").Value = -200805398
$ws.Range("N132").Value = -39059.999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1414.3024
$ws.Range("I61").Value = 945.7714
$ws.Range("J61").Value = 3464.125
$ws.Range("K61").Value = 945.7714
$ws.Range("L61").Value = 3464.125
$ws.Range("M61").Value = -733.7714
$ws.Range("N61").Value = -3888.125
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H109").Value = 31500
$ws.Range("J109").Value = 31500
$ws.Range("L109").Value = 31500
$ws.Range("N109").Value = -34274
$ws.Range("H112").Value = 29071.428
$ws.Range("J112").Value = 29071.428
$ws.Range("L112").Value = 29071.428
$ws.Range("N112").Value = -32025.428
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H115").Value = 28445
$ws.Range("J115").Value = 28445
$ws.Range("L115").Value = 28445
$ws.Range("N115").Value = -31579
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H119").Value = 35698
$ws.Range("J119").Value = 35698
$ws.Range("L119").Value = 35698
$ws.Range("N119").Value = -45374
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H132").Value = 2369.3147
$ws.Range("I132").Value = 1779.674
$ws.Range("J132").Value = 5759.75
$ws.Range("K132").Value = 5339.022
$ws.Range("L132").Value = 17279.25
$ws.Range("M132").Value = -2809.022
$ws.Range("N132").Value = -22339.25
$ws.Range("H136").Value = 1414.3024
$ws.Range("I136").Value = 945.7714
$ws.Range("J136").Value = 3464.125
$ws.Range("K136").Value = 2837.3142
$ws.Range("L136").Value = 10392.375
$ws.Range("M136").Value = -287.3141999999998
$ws.Range("N136").Value = -15492.375
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2163.8823
$ws.Range("I134").Value = 1188.8276
$ws.Range("J134").Value = 7819.2
$ws.Range("K134").Value = 3566.4828
$ws.Range("L134").Value = 23457.6
$ws.Range("M134").Value = -1031.4828
$ws.Range("N134").Value = -28527.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1840.859
$ws.Range("I58").Value = 1527.7612
$ws.Range("J58").Value = 3747.9092
$ws.Range("K58").Value = 1527.7612
$ws.Range("L58").Value = 3747.9092
$ws.Range("M58").Value = -1324.7612
$ws.Range("N58").Value = -4153.9092
$ws.Range("H134").Value = 9499.467000000001
$ws.Range("I134").Value = 11406.4
$ws.Range("K134").Value = 34219.2
$ws.Range("M134").Value = -31684.2
$ws.Range("H136").Value = 1840.859
$ws.Range("I136").Value = 1527.7612
$ws.Range("J136").Value = 3747.9092
$ws.Range("K136").Value = 4583.2836
$ws.Range("L136").Value = 11243.7276
$ws.Range("M136").Value = -2033.2836
$ws.Range("N136").Value = -16343.7276
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 6825
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 6825
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 20475
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -27359
$ws.Range("H132").Value = 2301.1428
$ws.Range("I132").Value = 916.7646999999999
$ws.Range("J132").Value = 4440.636
$ws.Range("K132").Value = 8250.882299999999
$ws.Range("L132").Value = 39965.724
$ws.Range("M132").Value = -5720.882299999999
$ws.Range("N132").Value = -45025.724
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4289.1333
$ws.Range("I132").Value = 2702
$ws.Range("K132").Value = 8106
$ws.Range("M132").Value = -5576
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3110.9167
$ws.Range("I122").Value = 2726.842
$ws.Range("J122").Value = 4570.4
$ws.Range("K122").Value = 8180.526
$ws.Range("L122").Value = 13711.2
$ws.Range("M122").Value = -5730.526
$ws.Range("N122").Value = -18611.2
$ws.Range("H132").Value = 7894.1816
$ws.Range("I132").Value = 2918.8
$ws.Range("J132").Value = 9357.529
$ws.Range("K132").Value = 8756.400000000001
$ws.Range("L132").Value = 28072.587
$ws.Range("M132").Value = -6226.400000000001
$ws.Range("N132").Value = -33132.587
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2549.9048
$ws.Range("I122").Value = 1796
$ws.Range("J122").Value = 4231.6924
$ws.Range("K122").Value = 5388
$ws.Range("L122").Value = 12695.0772
$ws.Range("M122").Value = -2938
$ws.Range("N122").Value = -17595.0772
$ws.Range("H126").Value = 275211.84
$ws.Range("I126").Value = 1876.64
$ws.Range("J126").Value = 763310.4399999999
$ws.Range("K126").Value = 5629.92
$ws.Range("L126").Value = 2289931.32
$ws.Range("M126").Value = -3159.92
$ws.Range("N126").Value = -2294871.32
$ws.Range("H136").Value = 3982.5454
$ws.Range("I136").Value = 856.4545000000001
$ws.Range("J136").Value = 7108.636
$ws.Range("K136").Value = 2569.3635
$ws.Range("L136").Value = 21325.908
$ws.Range("M136").Value = -19.36350000000039
$ws.Range("N136").Value = -26425.908
